{"js": "// Office.js (Word JavaScript API) script.\n// Implements:\n//   1. Append a new run \" (opcional)\" after \"Tela de senha;\".\n//   2. Insert a new paragraph with the \"Uma empresa de vans...\" text right\n//      before the trailing (bookmark) paragraph.\n//   3. Insert a new run \"Requisitos: \" at the very start of that trailing\n//      (bookmark) paragraph.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\n// --- Change 1: \"Tela de senha;\" -> add \" (opcional)\" run at paragraph end ---\nconst senhaPara = items.find((p) => p.text.trim() === \"Tela de senha;\");\nif (!senhaPara) {\n  throw new Error('Paragraph \"Tela de senha;\" not found.');\n}\nconst senhaEndRange = senhaPara.getRange(Word.RangeLocation.end);\n// Insert as raw OOXML so it lands as its own <w:r> (matching the target\n// diff) instead of being merged into the preceding run's text.\nconst opcionalOoxml =\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData>' +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  '<w:body><w:p><w:r><w:t xml:space=\"preserve\"> (opcional)</w:t></w:r></w:p></w:body>' +\n  '</w:document>' +\n  '</pkg:xmlData></pkg:part></pkg:package>';\nsenhaEndRange.insertOoxml(opcionalOoxml, Word.InsertLocation.end);\n\n// --- Change 2 & 3: new content right before the final (bookmark) paragraph ---\n// The trailing paragraph is the very last paragraph of the body (it holds\n// the _GoBack bookmark and, before this edit, no text runs).\nconst lastPara = items[items.length - 1];\n\nlastPara.insertParagraph(\n  \"Uma empresa de vans escolares no qual est\u00e1 em expans\u00e3o, solicitou um desenvolvimento de um programa no qual precisa armazenar os dados de suas vans, motoristas, alunos e as escolas de destino qual as crian\u00e7as pertencem e suas vans fazem o trajeto.    \",\n  Word.InsertLocation.before\n);\n\n// Add the \"Requisitos: \" run at the very beginning of the (still last)\n// bookmark paragraph, ahead of the bookmark.\nlastPara.insertText(\"Requisitos: \", Word.InsertLocation.start);\n\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) script.\n# Implements:\n#   1. Append \" (opcional)\" right after \"Tela de senha;\".\n#   2. Insert a new paragraph with the \"Uma empresa de vans...\" text right\n#      before the trailing (bookmark) paragraph.\n#   3. Insert \"Requisitos: \" at the very start of that trailing (bookmark)\n#      paragraph, ahead of the _GoBack bookmark.\n\n$d = $word.ActiveDocument\n\n# --- Change 1: \"Tela de senha;\" -> append \" (opcional)\" -------------------\n$rng = $d.Content\n$found = $rng.Find.Execute(\"Tela de senha;\")\nif ($found) {\n    $rng.Collapse(0)   # wdCollapseEnd\n    $rng.InsertAfter(\" (opcional)\")\n}\n\n# --- Change 2: new paragraph right before the trailing bookmark paragraph -\n$lastPara = $d.Paragraphs($d.Paragraphs.Count)\n$lastPara.Range.InsertParagraphBefore()\n\n$newPara = $d.Paragraphs($d.Paragraphs.Count - 1)\n$newPara.Range.Text = \"Uma empresa de vans escolares no qual est\u00e1 em expans\u00e3o, solicitou um desenvolvimento de um programa no qual precisa armazenar os dados de suas vans, motoristas, alunos e as escolas de destino qual as crian\u00e7as pertencem e suas vans fazem o trajeto.    \"\n\n# --- Change 3: \"Requisitos: \" at the start of the (still last) bookmark paragraph ---\n$bookmarkPara = $d.Paragraphs($d.Paragraphs.Count)\n$bookmarkPara.Range.InsertBefore(\"Requisitos: \")\n"}
